{"js": "// Task \"d\" (\u0413\u0430\u043b\u0430\u043a\u0442\u0438\u0447\u0435\u0441\u043a\u0430\u044f \u043d\u0430\u043b\u043e\u0433\u043e\u0432\u0430\u044f) limits fix:\n//   1) \"... n <= 10^5    m <= 100)\" -> \"... n <= 10^5    m < 90)\"\n//   2) \"... \u0437\u0430\u043f\u0440\u0435\u0449\u0451\u043d\u043d\u044b\u0445 \u0446\u0438\u0444\u0440 \u0432 \u043f\u043e\u0440\u044f\u0434\u043a\u0435 \u0432\u043e\u0437\u0440\u0430\u0441\u0442\u0430\u043d\u0438\u044f.\" -> \"... \u0437\u0430\u043f\u0440\u0435\u0449\u0451\u043d\u043d\u044b\u0445 \u0446\u0438\u0444\u0440.\"\n//\n// Both edits are done with Body.search() + Range.insertText(..., \"Replace\")\n// so the surrounding runs / formatting (bold, sz, lang, etc.) that are not\n// part of the changed text are left completely untouched.\n\n// --- 1) \"m <= 100)\" -> \"m < 90)\" --------------------------------------\n// Search only the part after \"m \" so the existing \"m\" run (lang=\"en-US\")\n// is left alone and only the numeric/comparison text is replaced.\nconst limitResults = context.document.body.search(\" <= 100)\", { matchCase: true });\nlimitResults.load(\"items,text\");\nawait context.sync();\n\nif (limitResults.items.length === 0) {\n  throw new Error('Could not find \" <= 100)\" to fix the m limit.');\n}\nlimitResults.items[0].insertText(\" < 90)\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 2) drop the now-obsolete \"\u0432 \u043f\u043e\u0440\u044f\u0434\u043a\u0435 \u0432\u043e\u0437\u0440\u0430\u0441\u0442\u0430\u043d\u0438\u044f\" requirement ------\nconst orderResults = context.document.body.search(\n  \"\u0437\u0430\u043f\u0440\u0435\u0449\u0451\u043d\u043d\u044b\u0445 \u0446\u0438\u0444\u0440 \u0432 \u043f\u043e\u0440\u044f\u0434\u043a\u0435 \u0432\u043e\u0437\u0440\u0430\u0441\u0442\u0430\u043d\u0438\u044f.\",\n  { matchCase: true }\n);\norderResults.load(\"items,text\");\nawait context.sync();\n\nif (orderResults.items.length === 0) {\n  throw new Error('Could not find the \"... \u0432 \u043f\u043e\u0440\u044f\u0434\u043a\u0435 \u0432\u043e\u0437\u0440\u0430\u0441\u0442\u0430\u043d\u0438\u044f.\" sentence tail.');\n}\norderResults.items[0].insertText(\"\u0437\u0430\u043f\u0440\u0435\u0449\u0451\u043d\u043d\u044b\u0445 \u0446\u0438\u0444\u0440.\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Task \"d\" (\u0413\u0430\u043b\u0430\u043a\u0442\u0438\u0447\u0435\u0441\u043a\u0430\u044f \u043d\u0430\u043b\u043e\u0433\u043e\u0432\u0430\u044f) limits fix:\n#   1) \"... n <= 10^5    m <= 100)\" -> \"... n <= 10^5    m < 90)\"\n#   2) \"... \u0437\u0430\u043f\u0440\u0435\u0449\u0451\u043d\u043d\u044b\u0445 \u0446\u0438\u0444\u0440 \u0432 \u043f\u043e\u0440\u044f\u0434\u043a\u0435 \u0432\u043e\u0437\u0440\u0430\u0441\u0442\u0430\u043d\u0438\u044f.\" -> \"... \u0437\u0430\u043f\u0440\u0435\u0449\u0451\u043d\u043d\u044b\u0445 \u0446\u0438\u0444\u0440.\"\n#\n# Both edits use Range.Find/Replacement so any run whose text is untouched\n# (e.g. the preceding \"m\" run, which carries lang=\"en-US\") keeps its own\n# formatting intact.\n\n$d = $word.ActiveDocument\n\n# --- 1) \"m <= 100)\" -> \"m < 90)\" ---------------------------------------\n# Search only the part after \"m \" so the existing \"m\" run is left alone\n# and only the numeric/comparison text is replaced.\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \" <= 100)\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \" < 90)\"\n$ok1 = $find1.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\nif (-not $ok1) {\n    throw 'edit.ps1: could not find \" <= 100)\" to fix the m limit.'\n}\n\n# --- 2) drop the now-obsolete \"\u0432 \u043f\u043e\u0440\u044f\u0434\u043a\u0435 \u0432\u043e\u0437\u0440\u0430\u0441\u0442\u0430\u043d\u0438\u044f\" requirement ------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"\u0437\u0430\u043f\u0440\u0435\u0449\u0451\u043d\u043d\u044b\u0445 \u0446\u0438\u0444\u0440 \u0432 \u043f\u043e\u0440\u044f\u0434\u043a\u0435 \u0432\u043e\u0437\u0440\u0430\u0441\u0442\u0430\u043d\u0438\u044f.\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"\u0437\u0430\u043f\u0440\u0435\u0449\u0451\u043d\u043d\u044b\u0445 \u0446\u0438\u0444\u0440.\"\n$ok2 = $find2.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\nif (-not $ok2) {\n    throw 'edit.ps1: could not find the \"... \u0432 \u043f\u043e\u0440\u044f\u0434\u043a\u0435 \u0432\u043e\u0437\u0440\u0430\u0441\u0442\u0430\u043d\u0438\u044f.\" sentence tail.'\n}\n"}
